# Apply manual_corrections.xlsx edits:
#  - Add a new manual correction row to unit_file for plant 6559, unit 7,
#    changing prime_mover from GT to CT.
#  - Remove rows for plants 63003 and 64850 from plant_file (both lat & lon
#    entries), keeping 54975 and 62262 (leading-zero plant ids removed so
#    they match the EIA-EPA crosswalk correctly).
#  - Update the active sheet / selections to match the final workbook state.

$wb = $excel.ActiveWorkbook

$wsEpa       = $wb.Worksheets.Item("epa_clean")
$wsEia       = $wb.Worksheets.Item("eia_clean")
$wsGenerator = $wb.Worksheets.Item("generator_file")
$wsUnit      = $wb.Worksheets.Item("unit_file")
$wsPlant     = $wb.Worksheets.Item("plant_file")

# --- unit_file: add new correction row for plant 6559 / unit 7 ---
$wsUnit.Range("A28").Value = "6559"
$wsUnit.Range("B28").Value = "7"
$wsUnit.Range("C28").Value = "GT"
$wsUnit.Range("D28").Value = "prime_mover"
$wsUnit.Range("E28").Value = "CT"

# --- plant_file: drop plants 63003 & 64850 (lat + lon rows) ---
$wsPlant.Rows("8:9").Delete() | Out-Null
$wsPlant.Rows("4:5").Delete() | Out-Null

# --- restore/update each sheet's last selection (order matters: the
#     last sheet activated below becomes the workbook's active tab) ---
$wsEia.Range("C1").Select() | Out-Null
$wsGenerator.Range("C8").Select() | Out-Null
$wsUnit.Range("A28").Select() | Out-Null
$wsPlant.Range("B11").Select() | Out-Null
